$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "weight" (column C) and "generating_number" (column D) values
# Row 2 (car): weight 5 -> 1, generating_number 12 -> 5
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 5

# Row 3 (bus): generating_number 2 -> 5
$ws.Range("D3").Value = 5

# Row 4 (truck): generating_number 2 -> 5
$ws.Range("D4").Value = 5

# Row 5 (motorcycle): generating_number 2 -> 5
$ws.Range("D5").Value = 5
